$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New conversation rows captured from the messaging export (rows 64-69)
$rows = @(
    @("2025-10-12 16:11:21", "Noah", 8450689526, "13052054965", "Hey man what’s up?"),
    @("2025-10-12 16:11:48", "Noah", 8450689526, "13052054965", "How is your day going?"),
    @("2025-10-12 16:12:08", "Noah", 8450689526, "13052054965", "What are you doing today?"),
    @("2025-10-12 16:12:22", "Noah", 8450689526, "13052054965", "This is a test message"),
    @("2025-10-12 16:15:58", "Noah", 8450689526, "13052054965", "One more test message"),
    @("2025-10-12 16:16:11", "Noah", 8450689526, "13052054965", "And another")
)

$startRow = 64
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]

    # Phone numbers are stored as text even though they look numeric -
    # a leading apostrophe forces Excel to keep the value as text.
    $dCell = $ws.Cells.Item($r, 4)
    $dCell.Value = "'" + $data[3]
    $dCell.Style = "Normal"

    $ws.Cells.Item($r, 5).Value = $data[4]

    # Media / Channel columns stay empty (but still present as text cells).
    $fCell = $ws.Cells.Item($r, 6)
    $fCell.Value = "'"
    $fCell.Style = "Normal"

    $gCell = $ws.Cells.Item($r, 7)
    $gCell.Value = "'"
    $gCell.Style = "Normal"
}
